$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AC, AD, AE
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match header styling of the existing header row (copy style from AB1)
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Restore the text values (paste special may have affected formatting only, but ensure values are correct)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill in the team record data for each player row (2-38): Wins=55, Losses=60, Ties=0
for ($r = 2; $r -le 38; $r++) {
    $ws.Cells.Item($r, 29).Value = 55
    $ws.Cells.Item($r, 30).Value = 60
    $ws.Cells.Item($r, 31).Value = 0
}
